$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final state of the transaction log (rows 2-15).
# Columns: E = Transaction Type, N = Payment Type, P = InternalComment, T = USD Amount
$rows = @(
    @{ Row = 2;  E = "Withdrawal"; N = "Crypto";       P = "ETH";         T = 261.57 },
    @{ Row = 3;  E = "Deposit";    N = "Crypto";       P = "ETH";         T = 341.28 },
    @{ Row = 4;  E = "Deposit";    N = "Crypto";       P = "ETH";         T = 596.35320000000002 },
    @{ Row = 5;  E = "Withdrawal"; N = "Wiretransfer";  P = "Anywires";    T = 1622.46 },
    @{ Row = 6;  E = "Withdrawal"; N = "Crypto";       P = "ETH";         T = 500.02510000000001 },
    @{ Row = 7;  E = "Withdrawal"; N = "Crypto";       P = "ETH";         T = 999.98659999999995 },
    @{ Row = 8;  E = "Deposit";    N = "Crypto";       P = "ETH";         T = 500.97559999999999 },
    @{ Row = 9;  E = "Withdrawal"; N = "Crypto";       P = "ETH";         T = 500 },
    @{ Row = 10; E = "Withdrawal"; N = "Crypto";       P = "ETH";         T = 540.65 },
    @{ Row = 11; E = "Withdrawal"; N = "Crypto";       P = "USDT ERC 20"; T = 1000.0771999999999 },
    @{ Row = 12; E = "Withdrawal"; N = "Crypto";       P = "ETH";         T = 1241.9655 },
    @{ Row = 13; E = "Deposit";    N = "Crypto";       P = "ETH";         T = 2496.0657999999999 },
    @{ Row = 14; E = "Deposit";    N = "Credit Card";  P = "Sipay";       T = 216.72 },
    @{ Row = 15; E = "Withdrawal"; N = "Credit Card";  P = "Startrading"; T = 250 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 14).Value = $r.N
    $ws.Cells.Item($r.Row, 16).Value = $r.P
    $ws.Cells.Item($r.Row, 20).Value = $r.T
}

# Highlight the newly-added top transaction (row 2, columns E:T) in red.
$ws.Range("E2:T2").Interior.Color = 255

# Move the active selection to G20, matching the saved workbook state.
$ws.Range("G20").Select()
